$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Empty (FPC Pascal)
$ws.Range("D3").Value = "33,792 bytes"
$ws.Range("E3").Value = "37,376 bytes"
$ws.Range("F3").Value = "36,352 bytes"

# Row 4 - Borland C
$ws.Range("D4").Value = "52,224 bytes"
$ws.Range("E4").Value = "66,560 bytes"
$ws.Range("F4").Value = "66,048 bytes"

# Row 5 - Borland C++
$ws.Range("D5").Value = "47,104 bytes"
$ws.Range("E5").Value = "149,504 bytes"
$ws.Range("F5").Value = "148,480 bytes"

# Row 6 - Dev-C++ .c files
$ws.Range("D6").Value = "331,079 bytes"
$ws.Range("E6").Value = "332307 bytes"
$ws.Range("F6").Value = "331,446 bytes"

# Row 7 - Dev-C++ .cpp files
$ws.Range("D7").Value = "331,090 bytes"
$ws.Range("E7").Value = "3,137,623 bytes"
$ws.Range("F7").Value = "3,137,114 bytes"

# Row 8 - Java
$ws.Range("D8").Value = "257 bytes"
$ws.Range("E8").Value = "907 bytes"
$ws.Range("F8").Value = "603 bytes"

# Row 9 - C#
$ws.Range("D9").Value = "35,328 bytes"
$ws.Range("E9").Value = "35,840 bytes"
$ws.Range("F9").Value = "35,328 bytes"

# Row 10 - Parva
$ws.Range("D10").Value = "22 bytes"
$ws.Range("E10").Value = "532 bytes"
$ws.Range("F10").Value = "220 bytes"

# Answers to the questions
$ws.Range("D13").Value = "Because they contain instructions to be used by the loader"
$ws.Range("D15").Value = "Pascal"
$ws.Range("D17").Value = "C#"
$ws.Range("D19").Value = "C#"
$ws.Range("D21").Value = "Dev-C++.cpp files because C++ IOstream relies on templates which create more inline conde"

# Update the saved selection
$ws.Range("G17").Select()
